$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.936.02"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "1.632.72"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'214.71"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.256"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'19.69"
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "1.856.95"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "1.625.42"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "'62.83"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "25.914.23"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D20").Value = "'193.56"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").Value = "'4.39"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").Value = "'9.96"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").Value = "'6.28"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").Value = "'142.26"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").Value = "'6.87"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "'15.48"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "'3.23"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("D36").Value = "'0.902"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "1.135.79"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").Value = "'2.47"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "'99.26"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("D44").Value = "'5.45"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "1.766.28"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "'56.20"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("E48").Value = "  +3.56%  "
$ws.Range("E49").Value = "  +1.92%  "
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").Value = "'7.61"
$ws.Range("E51").Value = "  +2.39%  "
